$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.362.29"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.07%  "
$ws.Range("D3").Value = "'1.825.59"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.43%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "'314.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.51%  "
$ws.Range("D6").Value = "'0.9998"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").Value = "'0.4673"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.18%  "
$ws.Range("D8").Value = "'0.3792"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.56%  "
$ws.Range("D9").Value = "'0.07428"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.56%  "
$ws.Range("D10").Value = "'0.8756"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.43%  "
$ws.Range("D11").Value = "'20.77"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.85%  "
$ws.Range("D12").Value = "'1.826.63"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.20%  "
$ws.Range("D13").Value = "'6.693"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.28%  "
$ws.Range("D14").Value = "'5.431"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.76%  "
$ws.Range("D15").Value = "'93.10"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.06%  "
$ws.Range("D16").Value = "'0.07090"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.17%  "
$ws.Range("E17").Value = "  -0.11%  "
$ws.Range("D18").Value = "'0.000008795"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.09%  "
$ws.Range("D19").Value = "'0.9997"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("D20").Value = "'15.03"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.34%  "
$ws.Range("D21").Value = "'27.358.86"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.93%  "
$ws.Range("D22").Value = "'5.323"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.49%  "
$ws.Range("E23").Value = "  +1.61%  "
$ws.Range("D24").Value = "'2.053.70"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.40%  "
$ws.Range("D25").Value = "'1.940"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.46%  "
$ws.Range("D26").Value = "'151.15"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.48%  "
$ws.Range("E27").Value = "  +3.53%  "
$ws.Range("D28").Value = "'18.67"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.21%  "
$ws.Range("D29").Value = "'5.335"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.82%  "
$ws.Range("D30").Value = "'117.31"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.67%  "
$ws.Range("D31").Value = "'0.08967"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.94%  "
$ws.Range("D32").Value = "'0.7911"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.22%  "
$ws.Range("D33").Value = "'1.195"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.60%  "
$ws.Range("E34").Value = "  +2.09%  "
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("D36").Value = "'0.9994"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.05%  "
$ws.Range("D37").Value = "'1.101"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.36%  "
$ws.Range("D38").Value = "'0.01976"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.85%  "
$ws.Range("D39").Value = "'0.05249"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.25%  "
$ws.Range("D40").Value = "'7.295"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.96%  "
$ws.Range("D41").Value = "'0.5342"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.68%  "
$ws.Range("B42").Value = "MXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D42").Value = "'2.896"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.46%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "'2.363"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +20.20%  "
$ws.Range("D44").Value = "'0.1704"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.31%  "
$ws.Range("E45").Value = "  +2.65%  "
$ws.Range("D46").Value = "'0.5093"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.60%  "
$ws.Range("D47").Value = "'10.59"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.40%  "
$ws.Range("D48").Value = "'105.53"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.17%  "
$ws.Range("D49").Value = "'1.681"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.22%  "
$ws.Range("D50").Value = "'0.9994"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.02%  "
$ws.Range("E51").Value = "  +1.17%  "
